$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: the date/time value was refreshed with new scrape data.
$ws.Range("A13").Value = 45813.39350918982

# New row 14: latest price entry appended by the scraper.
$ws.Range("A14").NumberFormat = $ws.Range("A13").NumberFormat
$ws.Range("A14").Value = 45814.39344563471
$ws.Range("B14").Value = "EVOWHEY PROTEIN"
$ws.Range("C14").Value = "2Kg"
$ws.Range("D14").Value = "37,90€"
